# Update the cryptocurrency price/volume table (cryptos.xlsx) with the
# latest scraped values from the Sat Apr 15 09:20:42 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.467.66'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '2.104.40'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5225'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4500'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.61'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +16.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08947'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.76%  '
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.41'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.26%  '
$ws.Range('D13').Value = '2.096.89'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.752'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.757'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.13%  '
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06607'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.298'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('D23').Value = '30.509.32'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.347'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.26%  '
$ws.Range('D26').Value = '2.342.11'
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.592'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '163.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.41%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  +2.51%  '
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.677'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +8.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.155'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.946'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.37'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.43%  '
$ws.Range('E37').Value = '  -1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06778'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.487'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2276'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6914'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.255'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.308'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.42%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '14.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6363'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('E48').Value = '  -0.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.245'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.222'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '82.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.23%  '
